# "changed moves in importer"
# Update the "Move Sheet" with real move data for Yoshi (Super Smash Bros
# Melee), replacing the placeholder scaffolding that used to live in
# column H and rows 6-14, and adding Type/Direction columns (C, D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Move Sheet")

# Drop the old "Name"/"Egg" scaffolding column (H) entirely.
$ws.Columns.Item(8).Delete()

# Drop the leftover numbered placeholder rows (6-14) below the real data.
$ws.Range("A6:A14").EntireRow.Delete()

# Data rows - Yoshi moves from Super Smash Bros Melee.
# Shared-string insertion order matters for byte-identical output, so we
# write column-by-column (A, B, C, D, ...) across all rows, matching how
# the importer tool produced the original file.
$ws.Range("A2").Value = "Yoshi"
$ws.Range("A3").Value = "Yoshi"
$ws.Range("A4").Value = "Yoshi"
$ws.Range("A5").Value = "Yoshi"

$ws.Range("B2").Value = "Super Smash Bros Melee"
$ws.Range("B3").Value = "Super Smash Bros Melee"
$ws.Range("B4").Value = "Super Smash Bros Melee"
$ws.Range("B5").Value = "Super Smash Bros Melee"

$ws.Range("C2").Value = "smash"
$ws.Range("C3").Value = "special"
$ws.Range("C4").Value = "strong"
$ws.Range("C5").Value = "smash"

$ws.Range("D2").Value = "up"
$ws.Range("D3").Value = "down"
$ws.Range("D4").Value = "up"
$ws.Range("D5").Value = "down"

$ws.Range("E2").Value = 10
$ws.Range("E3").Value = 12
$ws.Range("E4").Value = 9
$ws.Range("E5").Value = 2

$ws.Range("F2").Value = 15
$ws.Range("F3").Value = 13
$ws.Range("F4").Value = 11
$ws.Range("F5").Value = 5

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 4

# Column widths for the two new columns used by Type / Direction.
$ws.Columns.Item(3).ColumnWidth = 15.33203125
$ws.Columns.Item(4).ColumnWidth = 11.6640625

# Selection moves to D9 as in the saved file.
$ws.Range("D9").Select()

# Scroll position on "Game Data" (active tab) shifted too.
$wb.Worksheets.Item("Game Data").Activate()
$excel.ActiveWindow.ScrollRow = 2
